# Updated cryptos list on Tue Jul  4 21:43:48 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for existing rows, and
# updates rows 18-21 (coin name/link/price/volume) to reflect a re-ranked
# ShibaInu/Avalanche/Dai/WrappedliquidstakedEther2.0 ordering.
#
# Numeric-looking Price strings are assigned with a leading apostrophe so
# Excel keeps them as literal text (e.g. "0.000007717" would otherwise be
# auto-converted to a number), then the cell style is reset back to
# "Normal" so the quote-prefix formatting doesn't linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'30.800.60"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.91%  "
$c = $ws.Range("D3")
$c.Value = "'1.941.51"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.Value = "'242.43"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "
$c = $ws.Range("D6")
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$c = $ws.Range("D7")
$c.Value = "'0.4893"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.25%  "
$c = $ws.Range("D8")
$c.Value = "'0.2958"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.28%  "
$c = $ws.Range("D9")
$c.Value = "'0.06892"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.09%  "
$c = $ws.Range("D10")
$c.Value = "'19.45"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.90%  "
$c = $ws.Range("D11")
$c.Value = "'106.43"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "
$c = $ws.Range("D12")
$c.Value = "'1.956.70"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.08%  "
$c = $ws.Range("D13")
$c.Value = "'0.07728"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "
$c = $ws.Range("D14")
$c.Value = "'5.346"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.08%  "
$c = $ws.Range("D15")
$c.Value = "'0.6987"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.63%  "
$c = $ws.Range("D16")
$c.Value = "'277.26"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.83%  "
$c = $ws.Range("D17")
$c.Value = "'30.797.25"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D18")
$c.Value = "'0.000007717"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D19")
$c.Value = "'13.12"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D20")
$c.Value = "'0.9999"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D21")
$c.Value = "'2.193.09"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$c = $ws.Range("D22")
$c.Value = "'5.466"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("E23").Value = "  -0.36%  "
$c = $ws.Range("D24")
$c.Value = "'6.508"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("E25").Value = "  -2.19%  "
$c = $ws.Range("D26")
$c.Value = "'167.90"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "
$c = $ws.Range("D27")
$c.Value = "'19.65"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.36%  "
$c = $ws.Range("D28")
$c.Value = "'2.156"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.46%  "
$c = $ws.Range("D29")
$c.Value = "'0.1047"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "
$c = $ws.Range("D30")
$c.Value = "'1.392"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.31%  "
$c = $ws.Range("D31")
$c.Value = "'1.552"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.69%  "
$c = $ws.Range("D32")
$c.Value = "'4.551"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -4.62%  "
$c = $ws.Range("D33")
$c.Value = "'4.363"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.31%  "
$ws.Range("E34").Value = "  -2.91%  "
$c = $ws.Range("D35")
$c.Value = "'0.7517"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("E36").Value = "  -0.48%  "
$c = $ws.Range("D37")
$c.Value = "'0.9996"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
$c = $ws.Range("D38")
$c.Value = "'2.731"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "
$c = $ws.Range("D39")
$c.Value = "'0.01996"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.52%  "
$c = $ws.Range("D40")
$c.Value = "'2.662"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.66%  "
$c = $ws.Range("D41")
$c.Value = "'78.21"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +6.50%  "
$c = $ws.Range("D42")
$c.Value = "'6.503"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.43%  "
$c = $ws.Range("D43")
$c.Value = "'2.098"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.61%  "
$c = $ws.Range("D44")
$c.Value = "'0.9060"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.70%  "
$c = $ws.Range("D45")
$c.Value = "'108.02"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'0.4399"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.97%  "
$c = $ws.Range("D47")
$c.Value = "'0.9989"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "
$c = $ws.Range("D48")
$c.Value = "'7.738"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.73%  "
$c = $ws.Range("D49")
$c.Value = "'991.63"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.58%  "
$c = $ws.Range("D50")
$c.Value = "'0.1245"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.52%  "
$c = $ws.Range("D51")
$c.Value = "'9.263"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.42%  "
